$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# --- Title shape ("标题 1"): "English" -> "并发编程", lang en-US/zh-CN -> zh-CN/en-US ---
$title = $s.Shapes.Item(1)
$titleRange = $title.TextFrame.TextRange
$titleRange.Text = "并发编程"
$titleRange.LanguageID = "zh-CN"

# --- Subtitle shape ("副标题 2"): merge the "双" + "元音" runs into one run ---
$subtitle = $s.Shapes.Item(2)
$subRange = $subtitle.TextFrame.TextRange
# "--  " (4 chars) + "双元音" (3 chars) -> replace the last 3 chars (双元音) with the new phrase
$target = $subRange.Characters(5, 3)
$target.Text = "可见性、原子性、有序性"
